{"js": "// Replace the date line and all 25 three-digit-by-one-digit multiplication\n// answers in the table with their new values, via exact text search +\n// replace. Every source string is unique in the document, so a scoped\n// single-hit search/replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"2026-02-01 Sunday\", \"2026-02-02 Monday\"],\n  [\"347\u00d72=694\", \"965\u00d78=7720\"],\n  [\"295\u00d76=1770\", \"182\u00d78=1456\"],\n  [\"397\u00d73=1191\", \"615\u00d76=3690\"],\n  [\"446\u00d73=1338\", \"963\u00d78=7704\"],\n  [\"693\u00d74=2772\", \"423\u00d77=2961\"],\n  [\"280\u00d74=1120\", \"889\u00d73=2667\"],\n  [\"702\u00d73=2106\", \"196\u00d73=588\"],\n  [\"547\u00d74=2188\", \"490\u00d72=980\"],\n  [\"363\u00d73=1089\", \"448\u00d76=2688\"],\n  [\"390\u00d79=3510\", \"139\u00d73=417\"],\n  [\"579\u00d73=1737\", \"878\u00d72=1756\"],\n  [\"320\u00d78=2560\", \"756\u00d75=3780\"],\n  [\"550\u00d79=4950\", \"754\u00d79=6786\"],\n  [\"971\u00d79=8739\", \"323\u00d78=2584\"],\n  [\"582\u00d78=4656\", \"142\u00d79=1278\"],\n  [\"650\u00d75=3250\", \"532\u00d73=1596\"],\n  [\"976\u00d79=8784\", \"584\u00d79=5256\"],\n  [\"234\u00d77=1638\", \"957\u00d72=1914\"],\n  [\"540\u00d79=4860\", \"587\u00d74=2348\"],\n  [\"498\u00d77=3486\", \"318\u00d75=1590\"],\n  [\"356\u00d79=3204\", \"291\u00d72=582\"],\n  [\"531\u00d78=4248\", \"486\u00d76=2916\"],\n  [\"845\u00d73=2535\", \"486\u00d72=972\"],\n  [\"121\u00d75=605\", \"849\u00d79=7641\"],\n  [\"917\u00d72=1834\", \"511\u00d76=3066\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 25 three-digit-by-one-digit multiplication\n# answers in the table with their new values. Every source string is unique\n# in the document, so a plain Find/Replace (wdReplaceAll, which here only\n# ever matches a single occurrence) is safe and unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-02-01 Sunday\", \"2026-02-02 Monday\"),\n    @(\"347\u00d72=694\", \"965\u00d78=7720\"),\n    @(\"295\u00d76=1770\", \"182\u00d78=1456\"),\n    @(\"397\u00d73=1191\", \"615\u00d76=3690\"),\n    @(\"446\u00d73=1338\", \"963\u00d78=7704\"),\n    @(\"693\u00d74=2772\", \"423\u00d77=2961\"),\n    @(\"280\u00d74=1120\", \"889\u00d73=2667\"),\n    @(\"702\u00d73=2106\", \"196\u00d73=588\"),\n    @(\"547\u00d74=2188\", \"490\u00d72=980\"),\n    @(\"363\u00d73=1089\", \"448\u00d76=2688\"),\n    @(\"390\u00d79=3510\", \"139\u00d73=417\"),\n    @(\"579\u00d73=1737\", \"878\u00d72=1756\"),\n    @(\"320\u00d78=2560\", \"756\u00d75=3780\"),\n    @(\"550\u00d79=4950\", \"754\u00d79=6786\"),\n    @(\"971\u00d79=8739\", \"323\u00d78=2584\"),\n    @(\"582\u00d78=4656\", \"142\u00d79=1278\"),\n    @(\"650\u00d75=3250\", \"532\u00d73=1596\"),\n    @(\"976\u00d79=8784\", \"584\u00d79=5256\"),\n    @(\"234\u00d77=1638\", \"957\u00d72=1914\"),\n    @(\"540\u00d79=4860\", \"587\u00d74=2348\"),\n    @(\"498\u00d77=3486\", \"318\u00d75=1590\"),\n    @(\"356\u00d79=3204\", \"291\u00d72=582\"),\n    @(\"531\u00d78=4248\", \"486\u00d76=2916\"),\n    @(\"845\u00d73=2535\", \"486\u00d72=972\"),\n    @(\"121\u00d75=605\", \"849\u00d79=7641\"),\n    @(\"917\u00d72=1834\", \"511\u00d76=3066\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
